# Insert a new review row at row 3 (pushing the existing rows 3-11 down to
# 4-12), then populate the new row with the incoming review data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(3).Insert()

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Excelente atendimento "
$ws.Range("C3").Value = 45905.77066481482
$ws.Range("D3").Value = "ODgxNWYzNjgtNjZlYy00YjFiLWFiZGMtZGVkNWIxNGE4ZDY1OjU3MDE2"
